$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("A2").Value = "Última actualización: 08:45:36"
$ws.Range("A3").Value = "Total filas: 69"
$ws.Range("A23").Value = "06:17:28"
$ws.Range("C23").Value = "16_SANTA ANA"
$ws.Range("D23").Value = 64
$ws.Range("A24").Value = "05:57:13"
$ws.Range("C24").Value = "23_HERNANDEZ"
$ws.Range("D24").Value = 84
$ws.Range("C40").Value = "11_ETCHEVERRY"
$ws.Range("C41").Value = "15_ABASTO"
$ws.Range("A45").Value = "08:45:36"
$ws.Range("B45").Value = "08:46"
$ws.Range("C45").Value = "215C_EL PATO"
$ws.Range("D45").Value = 1
$ws.Range("A46").Value = "07:50:33"
$ws.Range("B46").Value = "08:49"
$ws.Range("D46").Value = 59
$ws.Range("A47").Value = "07:12:46"
$ws.Range("B47").Value = "08:51"
$ws.Range("D47").Value = 99
$ws.Range("A48").Value = "08:27:16"
$ws.Range("B48").Value = "08:52"
$ws.Range("C48").Value = "23_HERNANDEZ"
$ws.Range("D48").Value = 25
$ws.Range("A49").Value = "07:12:46"
$ws.Range("B49").Value = "08:53"
$ws.Range("D49").Value = 101
$ws.Range("A50").Value = "08:10:18"
$ws.Range("B50").Value = "08:54"
$ws.Range("C50").Value = "215B_EL PATO"
$ws.Range("D50").Value = 44
$ws.Range("A51").Value = "07:12:46"
$ws.Range("B51").Value = "08:57"
$ws.Range("D51").Value = 105
$ws.Range("A52").Value = "07:38:39"
$ws.Range("B52").Value = "08:58"
$ws.Range("C52").Value = "215A_EL PATO"
$ws.Range("D52").Value = 80
$ws.Range("A53").Value = "08:10:18"
$ws.Range("B53").Value = "09:05"
$ws.Range("C53").Value = "10_OLMOS"
$ws.Range("D53").Value = 55
$ws.Range("B54").Value = "09:06"
$ws.Range("C54").Value = "16_SANTA ANA"
$ws.Range("D54").Value = 88
$ws.Range("B55").Value = "09:16"
$ws.Range("C55").Value = "27_EL RETIRO"
$ws.Range("D55").Value = 98
$ws.Range("A56").Value = "07:38:39"
$ws.Range("C56").Value = "14_ABASTO"
$ws.Range("D56").Value = 99
$ws.Range("A57").Value = "08:27:16"
$ws.Range("B57").Value = "09:17"
$ws.Range("C57").Value = "27_EL RETIRO"
$ws.Range("D57").Value = 50
$ws.Range("A59").Value = "08:10:18"
$ws.Range("B59").Value = "09:18"
$ws.Range("C59").Value = "14_ABASTO"
$ws.Range("D59").Value = 68
$ws.Range("A60").Value = "07:38:39"
$ws.Range("B60").Value = "09:29"
$ws.Range("C60").Value = "10_OLMOS"
$ws.Range("D60").Value = 111
$ws.Range("B61").Value = "09:31"
$ws.Range("C61").Value = "16_SANTA ANA"
$ws.Range("D61").Value = 81
$ws.Range("A62").Value = "08:10:18"
$ws.Range("B62").Value = "09:36"
$ws.Range("C62").Value = "23_HERNANDEZ"
$ws.Range("D62").Value = 86
$ws.Range("A63").Value = "07:50:33"
$ws.Range("C63").Value = "15_ABASTO"
$ws.Range("D63").Value = 109
$ws.Range("A64").Value = "08:27:16"
$ws.Range("B64").Value = "09:39"
$ws.Range("C64").Value = "23_HERNANDEZ"
$ws.Range("D64").Value = 72
$ws.Range("A65").Value = "07:50:33"
$ws.Range("B65").Value = "09:41"
$ws.Range("D65").Value = 111
$ws.Range("A66").Value = "08:10:18"
$ws.Range("B66").Value = "09:42"
$ws.Range("C66").Value = "11_ETCHEVERRY"
$ws.Range("D66").Value = 92
$ws.Range("A67").Value = "07:50:33"
$ws.Range("B67").Value = "09:43"
$ws.Range("C67").Value = "16_P MOR-SANTA ANA"
$ws.Range("D67").Value = 113
$ws.Range("B68").Value = "09:53"
$ws.Range("C68").Value = "10_OLMOS"
$ws.Range("D68").Value = 103
$ws.Range("A69").Value = "08:10:18"
$ws.Range("B69").Value = "09:59"
$ws.Range("C69").Value = "215C_EL PATO"
$ws.Range("D69").Value = 109
$ws.Range("A70").Value = "08:37:25"
$ws.Range("B70").Value = "10:05"
$ws.Range("D70").Value = 88
$ws.Range("A71").Value = "08:10:18"
$ws.Range("B71").Value = "10:06"
$ws.Range("C71").Value = "14_ABASTO"
$ws.Range("D71").Value = 116
$ws.Range("A72").Value = "08:27:16"
$ws.Range("B72").Value = "10:13"
$ws.Range("C72").Value = "17X38_ROMERO"
$ws.Range("D72").Value = 106
$ws.Range("A73").Value = "08:37:25"
$ws.Range("B73").Value = "10:29"
$ws.Range("C73").Value = "15_ABASTO"
$ws.Range("D73").Value = 112
$ws.Range("E73").Value = "LP1912"
$ws.Range("A74").Value = "08:45:36"
$ws.Range("B74").Value = "10:44"
$ws.Range("C74").Value = "11X44_ETCHEVERRY"
$ws.Range("D74").Value = 119
$ws.Range("E74").Value = "LP1912"

# --- Sheet 2: LP1912-215 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").Value = "Última actualización: 08:45:36"
$ws.Range("A3").Value = "Total filas: 16"
$ws.Range("A16").Value = "08:45:36"
$ws.Range("B16").Value = "08:46"
$ws.Range("C16").Value = "215C_EL PATO"
$ws.Range("D16").Value = 1
$ws.Range("A17").Value = "07:12:46"
$ws.Range("B17").Value = "08:53"
$ws.Range("D17").Value = 101
$ws.Range("A18").Value = "08:10:18"
$ws.Range("B18").Value = "08:54"
$ws.Range("C18").Value = "215B_EL PATO"
$ws.Range("D18").Value = 44
$ws.Range("A19").Value = "07:12:46"
$ws.Range("B19").Value = "08:57"
$ws.Range("D19").Value = 105
$ws.Range("A20").Value = "07:38:39"
$ws.Range("B20").Value = "08:58"
$ws.Range("C20").Value = "215A_EL PATO"
$ws.Range("D20").Value = 80
$ws.Range("A21").Value = "08:10:18"
$ws.Range("B21").Value = "09:59"
$ws.Range("C21").Value = "215C_EL PATO"
$ws.Range("D21").Value = 109
$ws.Range("E21").Value = "LP1912"

# --- Sheet 3: 6203-6173 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Value = "Última actualización: 08:45:36"
$ws.Range("A3").Value = "Total filas: 12"
$ws.Range("A14").Value = "08:45:36"
$ws.Range("B14").Value = "10:12"
$ws.Range("D14").Value = 87
$ws.Range("A15").Value = "08:27:16"
$ws.Range("B15").Value = "10:13"
$ws.Range("C15").Value = "215C_LA PLATA"
$ws.Range("D15").Value = 106
$ws.Range("E15").Value = "L6203"
$ws.Range("A16").Value = "08:37:25"
$ws.Range("B16").Value = "10:30"
$ws.Range("C16").Value = "215B_LP-P MOR-1 Y 57"
$ws.Range("D16").Value = 113
$ws.Range("E16").Value = "L6173"
$ws.Range("A17").Value = "08:45:36"
$ws.Range("B17").Value = "10:31"
$ws.Range("C17").Value = "215A_LA PLATA"
$ws.Range("D17").Value = 106
$ws.Range("E17").Value = "L6173"
